# Rename the three header/footer logo pictures:
#   - BTec logo (first-page header)  : image2.jpg -> image1.jpg
#   - Pearson logo (default footer)  : image1.png -> image2.png
#   - Pearson logo (first-page footer): image1.png -> image2.png
#
# Word keeps these inline pictures on the document's single Section; the
# "default" Footer is index 1 (wdHeaderFooterPrimary), the "first page"
# Header/Footer is index 2 (wdHeaderFooterFirstPage).

$d = $word.ActiveDocument

$sec = $d.Sections.Item(1)

# --- First-page header: BTec_Logo-Orange -----------------------------
$btecHeader = $sec.Headers.Item(2)
if ($btecHeader.Exists -and $btecHeader.Range.InlineShapes.Count -ge 1) {
    $btecHeader.Range.InlineShapes.Item(1).Name = "image1.jpg"
}

# --- Default footer: PearsonLogo.png ----------------------------------
$pearsonFooterPrimary = $sec.Footers.Item(1)
if ($pearsonFooterPrimary.Exists -and $pearsonFooterPrimary.Range.InlineShapes.Count -ge 1) {
    $pearsonFooterPrimary.Range.InlineShapes.Item(1).Name = "image2.png"
}

# --- First-page footer: PearsonLogo.png -------------------------------
$pearsonFooterFirst = $sec.Footers.Item(2)
if ($pearsonFooterFirst.Exists -and $pearsonFooterFirst.Range.InlineShapes.Count -ge 1) {
    $pearsonFooterFirst.Range.InlineShapes.Item(1).Name = "image2.png"
}
